$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember column C's width so the newly inserted "Page_Object" column can match it
$descWidth = $ws.Columns.Item(3).ColumnWidth

# Insert a new column before the old column D ("Action_Keyword" shifts right to E,
# carrying its original width along with it)
$ws.Columns.Item(4).Insert()

# New column D matches the width used by column C (both "wide" description-style columns)
$ws.Columns.Item(4).ColumnWidth = $descWidth

# Header for the new "Page_Object" column
$ws.Range("D1").Value = "Page_Object"

# Object-repository keys for the rows that interact with a page element
$ws.Range("D4").Value = "txt_name"
$ws.Range("D5").Value = "txt_pass"
$ws.Range("D6").Value = "btn_login"

# The login-click step now uses the shorter "click" action keyword (now column E)
$ws.Range("E6").Value = "click"

# Match the saved selection from the edited workbook
$ws.Range("E6").Select()
